$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '64.271.01'
$ws.Range("E2").Value = '  +0.33%  '
# Row 3
$ws.Range("D3").Value = '2.771.92'
$ws.Range("E3").Value = '  -0.22%  '
# Row 4
$ws.Range("E4").Value = '  +0.10%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '577.96'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.60%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '160.19'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.02%  '
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.996'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.18%  '
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.602'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.45%  '
# Row 9
$ws.Range("E9").Value = '  -2.52%  '
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '5.88'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.46%  '
# Row 11
$ws.Range("E11").Value = '  +4.01%  '
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.387'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.15%  '
# Row 13
$ws.Range("D13").Value = '3.256.01'
$ws.Range("E13").Value = '  -0.28%  '
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.06'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.67%  '
# Row 15
$ws.Range("D15").Value = '63.914.43'
$ws.Range("E15").Value = '  -0.05%  '
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000153'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.19%  '
# Row 17
$ws.Range("D17").Value = '2.776.40'
$ws.Range("E17").Value = '  -0.22%  '
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.25'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.99%  '
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.86'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.85%  '
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '362.33'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.59%  '
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.68'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.42%  '
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.999'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.64%  '
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.531'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -7.64%  '
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '64.99'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.96%  '
# Row 25
$ws.Range("E25").Value = '  -2.09%  '
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.58'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.86%  '
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.28%  '
# Row 28
$ws.Range("D28").Value = '0.0₃0918'
$ws.Range("E28").Value = '  -3.67%  '
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.36'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.53%  '
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.98'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.70%  '
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.36'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +8.66%  '
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '168.46'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.43%  '
# Row 33
$ws.Range("B33").Value = 'NEARProtocol'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.01'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.60%  '
# Row 34
$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.52'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.64%  '
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '20.26'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.28%  '
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.998'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.01%  '
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.82'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.17%  '
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.01'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.09%  '
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '350.47'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.91%  '
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.33'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.25%  '
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.20'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.04%  '
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '39.22'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.16%  '
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '22.71'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.08%  '
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '21.61'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.13%  '
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0595'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.91%  '
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '137.65'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.02%  '
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.632'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.60%  '
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0255'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.32%  '
# Row 49
$ws.Range("E49").Value = '  -1.50%  '
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.997'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.16%  '
# Row 51
$ws.Range("B51").Value = 'Maker'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D51").Value = '2.137.63'
$ws.Range("E51").Value = '  -0.61%  '
